$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.361.82"
$ws.Range("E2").Value = "  -4.63%  "

$ws.Range("D3").Value = "1.565.03"
$ws.Range("E3").Value = "  -5.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.50%  "

$ws.Range("E7").Value = "  -2.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.08"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3401"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.169"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07662"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.09%  "

$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.051"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.934"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.56%  "

$ws.Range("D16").Value = "1.562.46"
$ws.Range("E16").Value = "  -5.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001130"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06740"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.249"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5291"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.28%  "

$ws.Range("D25").Value = "22.355.90"
$ws.Range("E25").Value = "  -4.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.388"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.794"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.54%  "

$ws.Range("E28").Value = "  -4.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.979"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.84%  "

$ws.Range("D32").Value = "1.738.56"
$ws.Range("E32").Value = "  -5.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.208"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.011"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.013"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08450"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.91%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02535"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2325"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.537"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06429"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.280"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6341"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5978"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.762"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.097"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.273"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.85%  "

